# BASILAN.xlsx update:
#  - TUBURAN: MAYOR race re-tallied with a new candidate (KALLAHAL, RAJIE)
#    inserted between the two existing rows; Over/Under-votes totals updated.
#  - CITY OF LAMITAN: MAYOR race candidate list/totals were a stray duplicate
#    of the VICE-MAYOR race; corrected to the real MAYOR results (which
#    happen to equal the VICE-MAYOR figures), removing the two extra rows.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ----------------------------------------------------------------------
# TUBURAN sheet
# ----------------------------------------------------------------------
$wsT = $wb.Worksheets.Item("TUBURAN")

# Insert a new row for the extra MAYOR candidate (pushes everything from
# the old row 354 down by one row).
$wsT.Rows.Item(354).Insert()

$wsT.Range("A352").Value = "KALLAHAL, GIMA (IND)"
Set-TextValue $wsT.Range("B352") "32"
Set-TextValue $wsT.Range("C352") "0.65 %"

$wsT.Range("A353").Value = "KALLAHAL, JHABER (PDPLBN)"
Set-TextValue $wsT.Range("B353") "4,809"
Set-TextValue $wsT.Range("C353") "98.80 %"

$wsT.Range("A354").Value = "KALLAHAL, RAJIE (IND)"
Set-TextValue $wsT.Range("B354") "26"
Set-TextValue $wsT.Range("C354") "0.53 %"

# MAYOR Over-votes/Under-votes/Valid votes/Votes-obtained totals row,
# now shifted down to row 357.
Set-TextValue $wsT.Range("A357") "42"
Set-TextValue $wsT.Range("B357") "907"
Set-TextValue $wsT.Range("D357") "4867"

# ----------------------------------------------------------------------
# CITY OF LAMITAN sheet
# ----------------------------------------------------------------------
$wsL = $wb.Worksheets.Item("CITY OF LAMITAN")

# The MAYOR candidate rows 354:355 (FURIGAY, ORIC and SAKKALAHUL,
# AL-RASHEED) are removed entirely, shifting everything below up by two.
$wsL.Range("A354:A355").EntireRow.Delete()

$wsL.Range("A352").Value = "FURIGAY, GEM (UNA)"
Set-TextValue $wsL.Range("B352") "33,576"
Set-TextValue $wsL.Range("C352") "96.66 %"

$wsL.Range("A353").Value = "INJANG, MICHAEL (IND)"
Set-TextValue $wsL.Range("B353") "1,158"
Set-TextValue $wsL.Range("C353") "3.33 %"

# MAYOR Over-votes/Under-votes/Valid votes/Votes-obtained totals row,
# now shifted up to row 356.
Set-TextValue $wsL.Range("A356") "18"
Set-TextValue $wsL.Range("B356") "6862"
Set-TextValue $wsL.Range("D356") "34734"
